# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet (the old N/O/P columns - "Late"/"Heading"/"Outstanding" - shift one
# place to the right, to O/P/Q), and make that sheet the active / selected
# tab (instead of "Input").

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N (shifts old N:P -> new O:Q).
$wsSchedule.Range("N1").EntireColumn.Insert()

# The newly inserted column keeps the width of its left neighbour (M).
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# Make "Repayment Schedule" the active sheet / selected tab, with P5
# selected (instead of "Input" being selected, and instead of the old
# A9:XFD9 selection on this sheet).
$wsSchedule.Activate()
$wsSchedule.Range("P5").Select()
